# Insert a new data row at row 120 (pushing existing rows 120..245 down to 121..246)
# and populate it with a new observation, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(120).EntireRow.Insert()

$ws.Cells.Item(120, 1).Value = 3
$ws.Cells.Item(120, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44539
$ws.Cells.Item(120, 5).Value = 5
$ws.Cells.Item(120, 6).Value = 100112009
$ws.Cells.Item(120, 7).Value = "Acelga"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 290
$ws.Cells.Item(120, 11).Value = 1800
$ws.Cells.Item(120, 12).Value = 2000
$ws.Cells.Item(120, 13).Value = 1890
$ws.Cells.Item(120, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(120, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(120, 16).Value = 315
$ws.Cells.Item(120, 17).Value = 6
$ws.Cells.Item(120, 18).Value = "Hortaliza"
